$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")
$ws.Activate()
$ws.Rows(6).Delete()
$wb.Save()
